$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2 and 3, and add new student rows 4-12
# Columns: A=id, B=email, C=name, D=password, E=isCampCommitee, F=isNewLogin, G=facultyId

$ws.Cells.Item(2, 1).Value = "C123413"
$ws.Cells.Item(2, 2).Value = "mary@mail.com"
$ws.Cells.Item(2, 3).Value = "Mary"
$ws.Cells.Item(2, 4).Value = "password"
$ws.Cells.Item(2, 5).Value = $false
$ws.Cells.Item(2, 6).Value = $true
$ws.Cells.Item(2, 7).Value = "98d0e59407f946b7aed49150ceba8627"

$ws.Cells.Item(3, 1).Value = "C133313"
$ws.Cells.Item(3, 2).Value = "caleb@mail.com"
$ws.Cells.Item(3, 3).Value = "Caleb"
$ws.Cells.Item(3, 4).Value = "password"
$ws.Cells.Item(3, 5).Value = $false
$ws.Cells.Item(3, 6).Value = $true
$ws.Cells.Item(3, 7).Value = "98d0e59407f946b7aed49150ceba8627"

$ws.Cells.Item(4, 1).Value = "C120513"
$ws.Cells.Item(4, 2).Value = "lopez@mail.com"
$ws.Cells.Item(4, 3).Value = "Lopez"
$ws.Cells.Item(4, 4).Value = "password"
$ws.Cells.Item(4, 5).Value = $false
$ws.Cells.Item(4, 6).Value = $true
$ws.Cells.Item(4, 7).Value = "98d0e59407f946b7aed49150ceba8627"

$ws.Cells.Item(5, 1).Value = "C121713"
$ws.Cells.Item(5, 2).Value = "winston@mail.com"
$ws.Cells.Item(5, 3).Value = "Winston"
$ws.Cells.Item(5, 4).Value = "password"
$ws.Cells.Item(5, 5).Value = $false
$ws.Cells.Item(5, 6).Value = $true
$ws.Cells.Item(5, 7).Value = "c9d8e441332d46bbb9655b8239c26e94"

$ws.Cells.Item(6, 1).Value = "C129013"
$ws.Cells.Item(6, 2).Value = "wick@mail.com"
$ws.Cells.Item(6, 3).Value = "Wick"
$ws.Cells.Item(6, 4).Value = "password"
$ws.Cells.Item(6, 5).Value = $false
$ws.Cells.Item(6, 6).Value = $true
$ws.Cells.Item(6, 7).Value = "c6a4beee0132472a99800fe0c310c731"

$ws.Cells.Item(7, 1).Value = "C127413"
$ws.Cells.Item(7, 2).Value = "will@mail.com"
$ws.Cells.Item(7, 3).Value = "Will"
$ws.Cells.Item(7, 4).Value = "password"
$ws.Cells.Item(7, 5).Value = $false
$ws.Cells.Item(7, 6).Value = $true
$ws.Cells.Item(7, 7).Value = "c9d8e441332d46bbb9655b8239c26e94"

$ws.Cells.Item(8, 1).Value = "C125413"
$ws.Cells.Item(8, 2).Value = "greg@mail.com"
$ws.Cells.Item(8, 3).Value = "Greg"
$ws.Cells.Item(8, 4).Value = "password"
$ws.Cells.Item(8, 5).Value = $false
$ws.Cells.Item(8, 6).Value = $true
$ws.Cells.Item(8, 7).Value = "98d0e59407f946b7aed49150ceba8627"

$ws.Cells.Item(9, 1).Value = "C128813"
$ws.Cells.Item(9, 2).Value = "henry@mail.com"
$ws.Cells.Item(9, 3).Value = "Henry"
$ws.Cells.Item(9, 4).Value = "password"
$ws.Cells.Item(9, 5).Value = $false
$ws.Cells.Item(9, 6).Value = $true
$ws.Cells.Item(9, 7).Value = "c6a4beee0132472a99800fe0c310c731"

$ws.Cells.Item(10, 1).Value = "C254513"
$ws.Cells.Item(10, 2).Value = "goh@mail.com"
$ws.Cells.Item(10, 3).Value = "Goh"
$ws.Cells.Item(10, 4).Value = "password"
$ws.Cells.Item(10, 5).Value = $false
$ws.Cells.Item(10, 6).Value = $true
$ws.Cells.Item(10, 7).Value = "98d0e59407f946b7aed49150ceba8627"

$ws.Cells.Item(11, 1).Value = "C028813"
$ws.Cells.Item(11, 2).Value = "hank@mail.com"
$ws.Cells.Item(11, 3).Value = "Hank"
$ws.Cells.Item(11, 4).Value = "password"
$ws.Cells.Item(11, 5).Value = $false
$ws.Cells.Item(11, 6).Value = $true
$ws.Cells.Item(11, 7).Value = "98d0e59407f946b7aed49150ceba8627"

$ws.Cells.Item(12, 1).Value = "C727413"
$ws.Cells.Item(12, 2).Value = "william@mail.com"
$ws.Cells.Item(12, 3).Value = "William"
$ws.Cells.Item(12, 4).Value = "password"
$ws.Cells.Item(12, 5).Value = $false
$ws.Cells.Item(12, 6).Value = $true
$ws.Cells.Item(12, 7).Value = "98d0e59407f946b7aed49150ceba8627"
